# Apply the "entity has changed, cvs added" update:
#  - SCORE (D2) and CLUSTER MUST HAVE SCORE (G2) recalculated
#  - CLUSTER MUST HAVE MATCH list (F2:F17) reordered/updated
#  - CLUSTER GOOD TO HAVE MATCH list (H2:H4) reordered

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scores -----------------------------------------------------------
# These look like numbers but must stay text cells (shared-string type),
# matching how the rest of the sheet stores its data. Using a leading
# apostrophe forces Excel to treat the input as text; resetting the
# style back to Normal afterwards drops the quote-prefix formatting that
# the apostrophe entry would otherwise leave behind.
$ws.Range("D2").Value = "'32.27"
$ws.Range("D2").Style = "Normal"

$ws.Range("G2").Value = "'33.33"
$ws.Range("G2").Style = "Normal"

# --- CLUSTER MUST HAVE MATCH (F2:F17) ----------------------------------
$ws.Range("F2").Value = "access : 1"
$ws.Range("F3").Value = "data architecture : 1"
$ws.Range("F4").Value = "sql : 4"
$ws.Range("F5").Value = "spark : 1"
$ws.Range("F6").Value = "python : 1"
$ws.Range("F7").Value = "data engineer : 1"
$ws.Range("F8").Value = "design : 4"
$ws.Range("F9").Value = "tools : 2"
$ws.Range("F10").Value = "aws : 2"
$ws.Range("F11").Value = "analysis : 6"
$ws.Range("F12").Value = "documentation : 2"
$ws.Range("F13").Value = "database : 1"
$ws.Range("F14").Value = "big data : 3"
$ws.Range("F15").Value = "engineer : 1"
$ws.Range("F16").Value = "analyze : 2"
$ws.Range("F17").Value = "hadoop : 1"

# --- CLUSTER GOOD TO HAVE MATCH (H2:H4) --------------------------------
$ws.Range("H2").Value = "aws : 2"
$ws.Range("H3").Value = "big data : 1"
$ws.Range("H4").Value = "management : 2"
